$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @(2,  45995,               1.154,               10.81),
  @(3,  45995.01041666666,   0.193,               3.597),
  @(4,  45995.02083333334,   4.634,               0.053),
  @(5,  45995.03125,         8.077,               0.774),
  @(6,  45995.04166666666,   0,                   21.36),
  @(7,  45995.05208333334,   0,                   20.499),
  @(8,  45995.0625,          4.299,               5.591),
  @(9,  45995.07291666666,   1.771,               1.467),
  @(10, 45995.08333333334,   0,                   29.264),
  @(11, 45995.09375,         0,                   36.298),
  @(12, 45995.10416666666,   0,                   15.642),
  @(13, 45995.11458333334,   0.956,               1.673),
  @(14, 45995.125,           0,                   18.749),
  @(15, 45995.13541666666,   0.092,               11.932),
  @(16, 45995.14583333334,   4.865,               0.367),
  @(17, 45995.15625,         0.33,                3.81),
  @(18, 45995.16666666666,   0.8100000000000001,  5.681),
  @(19, 45995.17708333334,   1.144,               6.498),
  @(20, 45995.1875,          7.342,               0),
  @(21, 45995.19791666666,   2.543,               8.712999999999999),
  @(22, 45995.20833333334,   0.881,               11.482),
  @(23, 45995.21875,         0.307,               1.253),
  @(24, 45995.22916666666,   5.354,               0.717),
  @(25, 45995.23958333334,   0.08,                6.546),
  @(26, 45995.25,            0,                   71.91200000000001),
  @(27, 45995.26041666666,   0,                   29.821),
  @(28, 45995.27083333334,   0,                   38.828),
  @(29, 45995.28125,         0,                   21.999),
  @(30, 45995.29166666666,   0.214,               8.041),
  @(31, 45995.30208333334,   6.71,                1.856),
  @(32, 45995.3125,          8.436999999999999,   0.403),
  @(33, 45995.32291666666,   4.627,               0.176),
  @(34, 45995.33333333334,   4.408,               7.875),
  @(35, 45995.34375,         0.004,               7.934),
  @(36, 45995.35416666666,   0.08400000000000001, 14.465),
  @(37, 45995.36458333334,   0,                   20.877),
  @(38, 45995.375,           0.067,               18.413),
  @(39, 45995.38541666666,   0,                   34.527),
  @(40, 45995.39583333334,   0,                   61.12),
  @(41, 45995.40625,         0,                   48.471),
  @(42, 45995.41666666666,   11.127,              0.379),
  @(43, 45995.42708333334,   0.008,               5.477),
  @(44, 45995.4375,          0,                   0)
)

foreach ($row in $data) {
  $r = $row[0]
  $ws.Cells.Item($r, 1).Value2 = $row[1]
  $ws.Cells.Item($r, 2).Value2 = $row[2]
  $ws.Cells.Item($r, 3).Value2 = $row[3]
}

# Apply the date/time style (same as used for existing column A cells) to the new row 44 A-cell
$ws.Range("A44").NumberFormat = $ws.Range("A43").NumberFormat
